$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new recipe data (Carotosa / carrot / 1 / facile / descrizione 5)
$ws.Range("A5").Value = "Carotosa"
$ws.Range("B5").Value = "carrot"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "facile"
$ws.Range("E5").Value = "descrizione 5"

# Update the active selection to match the edited range
$ws.Range("E4:E5").Select()
